$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Hassan Al Dhaheri"
$wsSummary.Range("B4").Value = 2313.86
$wsSummary.Range("B6").Value = 168025
$wsSummary.Range("B7").Value = 6014
$wsSummary.Range("B8").Value = 162011
$wsSummary.Range("B9").Value = 27.94

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("C2").Value = 165137
$wsAssets.Range("C3").Value = 2888
$wsAssets.Range("C4").Value = 168025

# ---------------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------------
# The "Auto Loans / Vehicle Loan 1" row is removed entirely; the
# "Credit Cards" row shifts up to row 2 (with updated figures) and the
# TOTAL LIABILITIES row shifts up to row 3 (with the new total).
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Rows.Item(2).Delete()

$wsLiabilities.Range("C2").Value = 6014
$wsLiabilities.Range("D2").Value = 301
$wsLiabilities.Range("E2").Value = 1

$wsLiabilities.Range("C3").Value = 6014
